# Insert a new price-report row at row 403 (weekly Choclo update for
# Feria Lagunitas de Puerto Montt). This pushes the existing rows 403-428
# down to 404-429 and fills the newly opened row 403 with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("403:403").Insert()

$ws.Range("A403").Value = 4
$ws.Range("B403").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C403").Value = "Los Lagos"
$ws.Range("D403").Value = 45223
$ws.Range("E403").Value = 10
$ws.Range("F403").Value = 100112024
$ws.Range("G403").Value = "Choclo"
$ws.Range("H403").Value = "Dulce o Americano"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 100
$ws.Range("K403").Value = 42000
$ws.Range("L403").Value = 42000
$ws.Range("M403").Value = 42000
$ws.Range("N403").Value = "$/malla 70 unidades"
$ws.Range("O403").Value = "Región de Arica y Parinacota"
$ws.Range("P403").Value = 600
$ws.Range("Q403").Value = 70
$ws.Range("R403").Value = "Hortaliza"
